# Apply row-level corrections to the Artfynd sheet per the source diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 111902034
$ws.Range("B2").Value = 90794
$ws.Range("D2").Value = "NT"
$ws.Range("E2").Value = 4362
$ws.Range("F2").Value = "Blå taggsvamp"
$ws.Range("G2").Value = "Hydnellum caeruleum"
$ws.Range("H2").Value = "(Hornem.) P.Karst."
$ws.Range("I2").NumberFormat = "@"
$ws.Range("I2").Value = "10"
$ws.Range("Q2").Value = 525039
$ws.Range("R2").Value = 6867407
$ws.Range("S2").Value = 25

# Row 3
$ws.Range("A3").Value = 111902031
$ws.Range("B3").Value = 90794
$ws.Range("E3").Value = 4362
$ws.Range("F3").Value = "Blå taggsvamp"
$ws.Range("G3").Value = "Hydnellum caeruleum"
$ws.Range("H3").Value = "(Hornem.) P.Karst."
$ws.Range("I3").NumberFormat = "@"
$ws.Range("I3").Value = "2"
$ws.Range("Q3").Value = 524990
$ws.Range("R3").Value = 6867386
$ws.Range("S3").Value = 25

# Row 4
$ws.Range("A4").Value = 111902029
$ws.Range("B4").Value = 88166
$ws.Range("D4").Value = "VU"
$ws.Range("E4").Value = 6276
$ws.Range("F4").Value = "Goliatmusseron"
$ws.Range("G4").Value = "Tricholoma matsutake"
$ws.Range("H4").Value = "(S.Ito & S.Imai) Singer"
$ws.Range("I4").NumberFormat = "@"
$ws.Range("I4").Value = "4"
$ws.Range("Q4").Value = 524972
$ws.Range("R4").Value = 6867342
$ws.Range("S4").Value = 5

# Row 5
$ws.Range("A5").Value = 111902040
$ws.Range("B5").Value = 90434
$ws.Range("E5").Value = 4745
$ws.Range("F5").Value = "Tallriska"
$ws.Range("G5").Value = "Lactarius musteus"
$ws.Range("H5").Value = "Fr."
$ws.Range("I5").NumberFormat = "@"
$ws.Range("I5").Value = "1"
$ws.Range("Q5").Value = 524891
$ws.Range("R5").Value = 6866840
$ws.Range("S5").Value = 10

# Row 6
$ws.Range("B6").Value = 88166

# Row 7
$ws.Range("A7").Value = 111902039
$ws.Range("B7").Value = 90816
$ws.Range("E7").Value = 2059
$ws.Range("F7").Value = "Skrovlig taggsvamp"
$ws.Range("G7").Value = "Hydnellum scabrosum"
$ws.Range("H7").Value = "(Fr.) E.Larss., K.H.Larss. & Kõljalg"
$ws.Range("I7").NumberFormat = "@"
$ws.Range("I7").Value = "5"
$ws.Range("Q7").Value = 524868
$ws.Range("R7").Value = 6867460
$ws.Range("S7").Value = 5

# Row 8
$ws.Range("A8").Value = 111902026
$ws.Range("B8").Value = 90816
$ws.Range("D8").Value = "NT"
$ws.Range("E8").Value = 2059
$ws.Range("F8").Value = "Skrovlig taggsvamp"
$ws.Range("G8").Value = "Hydnellum scabrosum"
$ws.Range("H8").Value = "(Fr.) E.Larss., K.H.Larss. & Kõljalg"
$ws.Range("Q8").Value = 524951
$ws.Range("R8").Value = 6867324

# Row 9
$ws.Range("A9").Value = 111902038
$ws.Range("B9").Value = 90800
$ws.Range("D9").Value = "LC"
$ws.Range("E9").Value = 4364
$ws.Range("F9").Value = "Dropptaggsvamp"
$ws.Range("G9").Value = "Hydnellum ferrugineum"
$ws.Range("H9").Value = "(Fr.:Fr.) P. Karst."
$ws.Range("Q9").Value = 524893
$ws.Range("R9").Value = 6867499

# Row 10
$ws.Range("A10").Value = 111902035
$ws.Range("B10").Value = 90792
$ws.Range("D10").Value = "NT"
$ws.Range("E10").Value = 4361
$ws.Range("F10").Value = "Orange taggsvamp"
$ws.Range("G10").Value = "Hydnellum aurantiacum"
$ws.Range("H10").Value = "(Batsch:Fr.) P.Karst."
$ws.Range("I10").NumberFormat = "@"
$ws.Range("I10").Value = "3"
$ws.Range("Q10").Value = 525047
$ws.Range("R10").Value = 6867385
$ws.Range("S10").Value = 25

# Row 11
$ws.Range("A11").Value = 111902030
$ws.Range("B11").Value = 88166
$ws.Range("D11").Value = "VU"
$ws.Range("E11").Value = 6276
$ws.Range("F11").Value = "Goliatmusseron"
$ws.Range("G11").Value = "Tricholoma matsutake"
$ws.Range("H11").Value = "(S.Ito & S.Imai) Singer"
$ws.Range("I11").NumberFormat = "@"
$ws.Range("I11").Value = "6"
$ws.Range("Q11").Value = 524971
$ws.Range("R11").Value = 6867379
$ws.Range("S11").Value = 5

# Row 12
$ws.Range("A12").Value = 111902028
$ws.Range("B12").Value = 90800
$ws.Range("D12").Value = "LC"
$ws.Range("E12").Value = 4364
$ws.Range("F12").Value = "Dropptaggsvamp"
$ws.Range("G12").Value = "Hydnellum ferrugineum"
$ws.Range("H12").Value = "(Fr.:Fr.) P. Karst."
$ws.Range("I12").NumberFormat = "@"
$ws.Range("I12").Value = "1"
$ws.Range("Q12").Value = 524954
$ws.Range("R12").Value = 6867304

# Row 13
$ws.Range("A13").Value = 111902027
$ws.Range("B13").Value = 90794
$ws.Range("I13").NumberFormat = "@"
$ws.Range("I13").Value = "5"
$ws.Range("Q13").Value = 524937
$ws.Range("R13").Value = 6867322

# Row 14
$ws.Range("A14").Value = 111902033
$ws.Range("B14").Value = 90434
$ws.Range("E14").Value = 4745
$ws.Range("F14").Value = "Tallriska"
$ws.Range("G14").Value = "Lactarius musteus"
$ws.Range("H14").Value = "Fr."
$ws.Range("I14").NumberFormat = "@"
$ws.Range("I14").Value = "1"
$ws.Range("Q14").Value = 525027
$ws.Range("R14").Value = 6867370
$ws.Range("S14").Value = 10

# Row 15
$ws.Range("A15").Value = 111902032
$ws.Range("B15").Value = 90792
$ws.Range("D15").Value = "NT"
$ws.Range("E15").Value = 4361
$ws.Range("F15").Value = "Orange taggsvamp"
$ws.Range("G15").Value = "Hydnellum aurantiacum"
$ws.Range("H15").Value = "(Batsch:Fr.) P.Karst."
$ws.Range("I15").NumberFormat = "@"
$ws.Range("I15").Value = "1"
$ws.Range("Q15").Value = 524989
$ws.Range("R15").Value = 6867384

# Row 16
$ws.Range("A16").Value = 111902037
$ws.Range("B16").Value = 90788
$ws.Range("D16").Value = "VU"
$ws.Range("E16").Value = 149
$ws.Range("F16").Value = "Tallgråticka"
$ws.Range("G16").Value = "Boletopsis grisea"
$ws.Range("H16").Value = "(Peck) Bondartsev & Singer"
$ws.Range("I16").NumberFormat = "@"
$ws.Range("I16").Value = "2"
$ws.Range("Q16").Value = 524869
$ws.Range("R16").Value = 6867441
